# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (column G) of the Session Analysis Results sheet.
# Cells currently reading "dnasr281@gmail.com, System" become
# "System, dnasr281@gmail.com". Cells with any other value (e.g. just
# "System" or just "dnasr281@gmail.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
